# Adds BOM, changes Pick and Place, and fixes component footprints
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pick Place for Buck_PCB")

# Component rows (Designator, Layer, Footprint, Rotation, Center X, Center Y)
# laid out in the exact column order they were entered: A, B, C, F, D, E.
$rows = @(
    @("C1", "TopLayer", "1210_-_Capacitor",     270, 15.748,  22.225),
    @("C2", "TopLayer", "0805_-_Capacitor",     90,  26.162,  12.664),
    @("C3", "TopLayer", "0805_-_Capacitor",     180, 19.05,   22.987),
    @("C4", "TopLayer", "0805_-_Capacitor",     180, 16.764,  12.954),
    @("C5", "TopLayer", "0805_-_Capacitor",     90,  15.875,  17.399),
    @("C6", "TopLayer", "0805_-_Capacitor",     90,  10.2498, 16.891),
    @("C7", "TopLayer", "0805_-_Capacitor",     180, 24.638,  22.987),
    @("C8", "TopLayer", "OCVZ0606_-_CAPACITOR", 270, 32.004,  9.779),
    @("D1", "TopLayer", "SMA_-_DIODE",          90,  28.3423, 19.812),
    @("L1", "TopLayer", "MSS1210_-_INDUCTOR",   0,   38.735,  21.971),
    @("P1", "TopLayer", "HDR1X2",               270, 12.192,  23.114),
    @("P2", "TopLayer", "HDR1X2",               270, 38.354,  11.049),
    @("R1", "TopLayer", "0603_-_Resistor",      90,  13.9208, 14.478),
    @("R2", "TopLayer", "0603_-_Resistor",      360, 12.954,  18.161),
    @("R3", "TopLayer", "0805_-_Resistor",      360, 20.574,  10.541),
    @("R4", "TopLayer", "0805_-_Resistor",      90,  12.2165, 14.478),
    @("U1", "TopLayer", "LM5005",               0,   21.0328, 17.78)
)

# Start with a clean sheet.
$ws.Cells.Clear()

# Header row: Designator, Layer, Footprint, Rotation first (Center X / Center Y
# headers are added last, after all the component data below).
$ws.Cells.Item(1, 1).Value = "Designator"
$ws.Cells.Item(1, 2).Value = "Layer"
$ws.Cells.Item(1, 3).Value = "Footprint"
$ws.Cells.Item(1, 6).Value = "Rotation"

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $r++
}

# Center X / Center Y headers, added last.
$ws.Cells.Item(1, 4).Value = "Center X"
$ws.Cells.Item(1, 5).Value = "Center Y"

# Adjust column widths to fit the new, longer footprint/center values.
$ws.Columns.Item(3).ColumnWidth = 21
$ws.Columns.Item(4).ColumnWidth = 11.5
$ws.Columns.Item(5).ColumnWidth = 11.5

# Set the active selection to G1 to match the saved file state.
$ws.Range("G1").Select() | Out-Null

Write-Output "Pick and place table updated."
